$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 57
$ws1.Range("F3").Value = 57
$ws1.Range("F4").Value = 952
$ws1.Range("F5").Value = 1246
$ws1.Range("F6").Value = 1716
$ws1.Range("F7").Value = 907
$ws1.Range("F8").Value = 566
$ws1.Range("F9").Value = 2507
$ws1.Range("F10").Value = 725
$ws1.Range("F11").Value = 571
$ws1.Range("F12").Value = 564
$ws1.Range("F13").Value = 21
$ws1.Range("F16").Value = 222
$ws1.Range("F17").Value = 514
$ws1.Range("F18").Value = 2111
$ws1.Range("F19").Value = 1227
$ws1.Range("F20").Value = 703
$ws1.Range("F22").Value = 2611
$ws1.Range("F25").Value = 517
$ws1.Range("F26").Value = 508
$ws1.Range("F27").Value = 291
$ws1.Range("F28").Value = 291
$ws1.Range("F30").Value = 1769
$ws1.Range("F33").Value = 519
$ws1.Range("F34").Value = 544
$ws1.Range("F36").Value = 4563
$ws1.Range("F37").Value = 133
$ws1.Range("F38").Value = 74

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 4196
$ws2.Range("F13").Value = 9
$ws2.Range("F15").Value = 319
$ws2.Range("F26").Value = 249
$ws2.Range("F28").Value = 255
$ws2.Range("F37").Value = 474

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1418
$ws3.Range("F6").Value = 517
$ws3.Range("F7").Value = 69
$ws3.Range("F8").Value = 186

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1418
$ws4.Range("F6").Value = 57
$ws4.Range("F8").Value = 57
$ws4.Range("F9").Value = 1246
$ws4.Range("F10").Value = 1716
$ws4.Range("F14").Value = 907
$ws4.Range("F15").Value = 566
$ws4.Range("F16").Value = 2507
$ws4.Range("F17").Value = 725
$ws4.Range("F18").Value = 571
$ws4.Range("F19").Value = 564
$ws4.Range("F20").Value = 21
$ws4.Range("F24").Value = 222
$ws4.Range("F25").Value = 9
$ws4.Range("F26").Value = 319
$ws4.Range("F28").Value = 514
$ws4.Range("F29").Value = 2111
$ws4.Range("F35").Value = 517
$ws4.Range("F36").Value = 186
$ws4.Range("F39").Value = 508
$ws4.Range("F40").Value = 291
$ws4.Range("F41").Value = 1769
$ws4.Range("F42").Value = 249
$ws4.Range("F43").Value = 519
$ws4.Range("F44").Value = 544
$ws4.Range("F46").Value = 4563
$ws4.Range("F47").Value = 133
$ws4.Range("F50").Value = 74
